$d = $word.ActiveDocument

# The target paragraph reads:
#   "En esta sección, nos vamos a enfocar en la creación de los datos del
#    proyecto, una parte crucial para el funcionamiento de la aplicación. "
# and should become:
#   "En esta sección se procederá a crear los datos del proyecto, una parte
#    crucial para el funcionamiento de la aplicación. "
# Only the opening clause changes; the trailing
# "los datos del proyecto, una parte crucial ... aplicación. " stays intact.

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Replacement.ClearFormatting()

$range.Find.Execute(
    "En esta sección, nos vamos a enfocar en la creación de los datos del proyecto",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "En esta sección se procederá a crear los datos del proyecto",
    2
)
